# Apply the "Bugfixes, mainly related to folders and experiments names" edit:
# add two new BlockName/Value rows (percent_reflecting_sfc=0.9 and R=0) at the
# bottom of Sheet1, and move the sheet's selection/scroll to reflect working
# near the newly appended rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# New row 29: percent_reflecting_sfc / 0.9
$ws.Cells.Item(29, 1).Value = "percent_reflecting_sfc"
$ws.Cells.Item(29, 2).Value = 0.9

# New row 30: R / 0
$ws.Cells.Item(30, 1).Value = "R"
$ws.Cells.Item(30, 2).Value = 0

# Match the final view state: scrolled down with B31 selected (first empty
# cell right below the data that was just added).
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("B31").Select()
